$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.730.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.236.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -8.76%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "177.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -13.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "511.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.592"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.234.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.609"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.129"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.766.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.77%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.245.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.57%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.115"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "62.744.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.935"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "365.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -11.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.27%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.47%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "626.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.104"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.23%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.386"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.40%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.918.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.122"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0638"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -15.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0383"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.24%  "
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -18.27%  "
